$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Update the "Recorded By" column (G) for rows where the list of recorders was reordered.
$ws.Cells.Item(2, 7).Value = 'Amira.Sobhy@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Cells.Item(3, 7).Value = 'eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Cells.Item(4, 7).Value = 'eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Cells.Item(5, 7).Value = 'Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(6, 7).Value = 'majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, mennatulla.medhat@med.asu.edu.eg'
$ws.Cells.Item(7, 7).Value = 'Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Cells.Item(9, 7).Value = 'Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Cells.Item(11, 7).Value = 'Safa.hany@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Cells.Item(12, 7).Value = 'amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg'
$ws.Cells.Item(13, 7).Value = 'esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Cells.Item(15, 7).Value = 'Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Cells.Item(17, 7).Value = 'esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Cells.Item(19, 7).Value = 'mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
$ws.Cells.Item(20, 7).Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Cells.Item(24, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg'
$ws.Cells.Item(25, 7).Value = 'Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg'
$ws.Cells.Item(28, 7).Value = 'maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Cells.Item(30, 7).Value = 'aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg'
